$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force text format on the Importe (H) column so that
# numeric-looking strings are stored as text, not auto-converted to numbers.
$importeRange = $ws.Range("H2:H112")
$importeRange.NumberFormat = "@"

$ws.Range("E49").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E77").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E50").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F50").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E78").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F78").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E51").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("H2").Value = "510.00"
$ws.Range("H3").Value = "270.00"
$ws.Range("H4").Value = "1013.00"
$ws.Range("H5").Value = "5074.14"
$ws.Range("H6").Value = "40000.00"
$ws.Range("H7").Value = "529.21"
$ws.Range("H8").Value = "1538.54"
$ws.Range("H9").Value = "1178.98"
$ws.Range("H10").Value = "188.00"
$ws.Range("H11").Value = "1496.00"
$ws.Range("H12").Value = "27790.00"
$ws.Range("H13").Value = "114719.08"
$ws.Range("H14").Value = "19911.04"
$ws.Range("H15").Value = "6917.32"
$ws.Range("H16").Value = "7133.61"
$ws.Range("H17").Value = "4650.00"
$ws.Range("H18").Value = "1086.00"
$ws.Range("H19").Value = "72.00"
$ws.Range("H20").Value = "662.00"
$ws.Range("H21").Value = "21520.00"
$ws.Range("H22").Value = "150.00"
$ws.Range("H72").Value = "150.00"
$ws.Range("H23").Value = "7407.50"
$ws.Range("H24").Value = "615.54"
$ws.Range("H25").Value = "27566.60"
$ws.Range("H26").Value = "1597.98"
$ws.Range("H27").Value = "11130.00"
$ws.Range("H28").Value = "6174.00"
$ws.Range("H29").Value = "25259.11"
$ws.Range("H30").Value = "1756.48"
$ws.Range("H31").Value = "200.75"
$ws.Range("H32").Value = "99092.50"
$ws.Range("H33").Value = "3071.00"
$ws.Range("H34").Value = "102.00"
$ws.Range("H35").Value = "1023.00"
$ws.Range("H36").Value = "52.00"
$ws.Range("H37").Value = "13.60"
$ws.Range("H38").Value = "412.00"
$ws.Range("H39").Value = "832.00"
$ws.Range("H40").Value = "5415.00"
$ws.Range("H41").Value = "452.00"
$ws.Range("H42").Value = "6350.00"
$ws.Range("H43").Value = "60.00"
$ws.Range("H44").Value = "1399.98"
$ws.Range("H45").Value = "28800.00"
$ws.Range("H46").Value = "1440.00"
$ws.Range("H47").Value = "4677.40"
$ws.Range("H48").Value = "264.00"
$ws.Range("H49").Value = "9919.25"
$ws.Range("H50").Value = "18639.87"
$ws.Range("H51").Value = "307.00"
$ws.Range("H52").Value = "1900.00"
$ws.Range("H53").Value = "4.20"
$ws.Range("H54").Value = "179500.00"
$ws.Range("H55").Value = "0.02"
$ws.Range("H56").Value = "103770.00"
$ws.Range("H57").Value = "148.58"
$ws.Range("H58").Value = "7052.00"
$ws.Range("H59").Value = "1517.00"
$ws.Range("H60").Value = "225.00"
$ws.Range("H61").Value = "4860.00"
$ws.Range("H62").Value = "130.20"
$ws.Range("H63").Value = "5740.00"
$ws.Range("H64").Value = "480.00"
$ws.Range("H65").Value = "39968.00"
$ws.Range("H66").Value = "18600.00"
$ws.Range("H67").Value = "980.00"
$ws.Range("H68").Value = "500.00"
$ws.Range("H69").Value = "657.28"
$ws.Range("H70").Value = "110.34"
$ws.Range("H71").Value = "428.00"
$ws.Range("H73").Value = "4500.00"
$ws.Range("H74").Value = "1260.00"
$ws.Range("H75").Value = "405.00"
$ws.Range("H76").Value = "456.00"
$ws.Range("H77").Value = "119.00"
$ws.Range("H78").Value = "250.87"
$ws.Range("H79").Value = "230.00"
$ws.Range("H80").Value = "1114.00"
$ws.Range("H81").Value = "673.04"
$ws.Range("H82").Value = "9075.00"
$ws.Range("H83").Value = "751.44"
$ws.Range("H84").Value = "2686.99"
$ws.Range("H85").Value = "2495.90"
$ws.Range("H86").Value = "19449.00"
$ws.Range("H87").Value = "94.00"
$ws.Range("H88").Value = "1599.64"
$ws.Range("H89").Value = "1625.08"
$ws.Range("H90").Value = "93.00"
$ws.Range("H91").Value = "1238.80"
$ws.Range("H92").Value = "540.00"
$ws.Range("H93").Value = "37.50"
$ws.Range("H94").Value = "525.00"
$ws.Range("H95").Value = "88.00"
$ws.Range("H96").Value = "23011.70"
$ws.Range("H97").Value = "2566.88"
$ws.Range("H98").Value = "42000.00"
$ws.Range("H99").Value = "350956.56"
$ws.Range("H100").Value = "2326.96"
$ws.Range("H101").Value = "380000.00"
$ws.Range("H106").Value = "380000.00"
$ws.Range("H102").Value = "34000.00"
$ws.Range("H103").Value = "358418.00"
$ws.Range("H104").Value = "243986.00"
$ws.Range("H105").Value = "279416.00"
$ws.Range("H107").Value = "360000.00"
$ws.Range("H108").Value = "199202.00"
$ws.Range("H109").Value = "44000.00"
$ws.Range("H110").Value = "284.75"
$ws.Range("H111").Value = "1228.00"
$ws.Range("H112").Value = "825000.00"

# Restore the original (default) style on the Importe column now that
# the values are committed as text, so no stray number-format style lingers.
$importeRange.Style = "Normal"

